# Fix: prevent hidden columns from being labeled upon detecting changes.
#
# The generator used to flag a row as "changed" (ÄNDERUNG, column L) even when
# the only "differences" it found were in columns that should have been
# skipped while diffing. This produced false positives for every row that:
#   - is the first row of a new "Datenelement" group (column B), and
#   - actually has identical content in the FV2210 block (C:K) and the
#     FV2304 block (N:V).
#
# For such rows the whole row must get the normal "group header" shading
# (the same direct formatting already used correctly elsewhere in the sheet,
# e.g. row 133) and column L must be cleared of the false "ÄNDERUNG" label.
# For continuation rows (same Datenelement as the row above) only column L
# needs to be reset - the rest of the row formatting was already correct.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Known-good reference row: a "new group" row with no real FV2210/FV2304
# differences, already carrying the correct group-header formatting
# (A=2 gray, B=3 bold+gray, C:K=2, L=4 gray+centered+empty, M:V=2).
$refRow = 133
$refFull = $ws.Range("A" + $refRow + ":V" + $refRow)
$refL = $ws.Range("L" + $refRow)

# Rows that are the first row of their Datenelement group and were wrongly
# flagged as changed -> restore full-row group-header formatting.
$fullRows = @(141,144,148,153,156,160,162,166,170,174,178,184,190,193,196)

# Continuation rows (same Datenelement as previous row) that were wrongly
# flagged as changed -> only column L needs to be reset.
$lOnlyRows = @(142,143,145,146,147,149,150,151,152,154,155,157,158,159,161,163,164,165,167,168,169,171,172,173,175,176,177,179,180,181,182,183,185,186,188,189,191,192,194,195,197,198)

foreach ($r in $fullRows) {
    $dst = $ws.Range("A" + $r + ":V" + $r)
    $refFull.Copy()
    $dst.PasteSpecial(-4122) # xlPasteFormats
    $ws.Range("L" + $r).ClearContents()
}

foreach ($r in $lOnlyRows) {
    $dstL = $ws.Range("L" + $r)
    $refL.Copy()
    $dstL.PasteSpecial(-4122) # xlPasteFormats
    $dstL.ClearContents()
}

$excel.CutCopyMode = 0
